$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.016.64"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "3.016.13"
$ws.Range("E3").Value = "  +3.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "615.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.16%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("E9").Value = "  +6.68%  "
$ws.Range("D10").Value = "3.013.01"
$ws.Range("E10").Value = "  +3.92%  "
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.19"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.54%  "
$ws.Range("D14").Value = "3.568.94"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("D16").Value = "75.928.20"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").Value = "3.007.36"
$ws.Range("E18").Value = "  +3.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.77%  "
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "3.172.43"
$ws.Range("E24").Value = "  +4.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.996"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.40%  "
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "490.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.57%  "
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.122"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.50%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "190.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.375"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("E42").Value = "  -4.62%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.774"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +18.94%  "
$ws.Range("E46").Value = "  +6.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "41.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("E49").Value = "  +8.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.589"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.21%  "
$ws.Range("E51").Value = "  +0.75%  "
